$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Row 19 <-> Row 20 swap (Uniswap/ShibaInu order swap) plus value updates
Set-TextValue $ws.Range("B19") "ShibaInu"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D19") "0.000007964"
Set-TextValue $ws.Range("E19") "  +5.89%  "

Set-TextValue $ws.Range("B20") "Uniswap"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D20") "5.827"
Set-TextValue $ws.Range("E20") "  +10.45%  "

# Remaining per-row Price(D) / Volume(E) updates
# Row 2
Set-TextValue $ws.Range("D2") "31.382.09"
Set-TextValue $ws.Range("E2") "  +3.24%  "
# Row 3
Set-TextValue $ws.Range("D3") "1.999.43"
Set-TextValue $ws.Range("E3") "  +6.79%  "
# Row 4
Set-TextValue $ws.Range("D4") "0.9987"
Set-TextValue $ws.Range("E4") "  -0.12%  "
# Row 5
Set-TextValue $ws.Range("D5") "0.8116"
Set-TextValue $ws.Range("E5") "  +72.25%  "
# Row 6
Set-TextValue $ws.Range("D6") "255.18"
Set-TextValue $ws.Range("E6") "  +4.50%  "
# Row 7
Set-TextValue $ws.Range("D7") "0.9977"
Set-TextValue $ws.Range("E7") "  -0.22%  "
# Row 8
Set-TextValue $ws.Range("D8") "0.3560"
Set-TextValue $ws.Range("E8") "  +23.31%  "
# Row 9
Set-TextValue $ws.Range("D9") "25.89"
Set-TextValue $ws.Range("E9") "  +17.66%  "
# Row 10
Set-TextValue $ws.Range("D10") "0.07026"
Set-TextValue $ws.Range("E10") "  +8.06%  "
# Row 11
Set-TextValue $ws.Range("D11") "0.8463"
Set-TextValue $ws.Range("E11") "  +16.79%  "
# Row 12
Set-TextValue $ws.Range("D12") "0.08140"
Set-TextValue $ws.Range("E12") "  +4.48%  "
# Row 13
Set-TextValue $ws.Range("D13") "101.50"
Set-TextValue $ws.Range("E13") "  +5.38%  "
# Row 14
Set-TextValue $ws.Range("D14") "1.996.19"
Set-TextValue $ws.Range("E14") "  +6.64%  "
# Row 15
Set-TextValue $ws.Range("D15") "5.528"
Set-TextValue $ws.Range("E15") "  +7.37%  "
# Row 16
Set-TextValue $ws.Range("D16") "273.01"
Set-TextValue $ws.Range("E16") "  -3.43%  "
# Row 17
Set-TextValue $ws.Range("D17") "31.385.28"
Set-TextValue $ws.Range("E17") "  +3.26%  "
# Row 18
Set-TextValue $ws.Range("E18") "  +7.72%  "
# Row 21
Set-TextValue $ws.Range("D21") "2.254.92"
Set-TextValue $ws.Range("E21") "  +6.81%  "
# Row 22
Set-TextValue $ws.Range("D22") "0.9971"
Set-TextValue $ws.Range("E22") "  -0.27%  "
# Row 23
Set-TextValue $ws.Range("D23") "0.9987"
Set-TextValue $ws.Range("E23") "  -0.11%  "
# Row 24
Set-TextValue $ws.Range("D24") "7.037"
Set-TextValue $ws.Range("E24") "  +12.57%  "
# Row 25
Set-TextValue $ws.Range("D25") "9.876"
Set-TextValue $ws.Range("E25") "  +8.84%  "
# Row 26
Set-TextValue $ws.Range("D26") "0.1524"
Set-TextValue $ws.Range("E26") "  +58.37%  "
# Row 27
Set-TextValue $ws.Range("D27") "164.44"
Set-TextValue $ws.Range("E27") "  +0.37%  "
# Row 28
Set-TextValue $ws.Range("D28") "20.15"
Set-TextValue $ws.Range("E28") "  +7.02%  "
# Row 29
Set-TextValue $ws.Range("D29") "2.274"
Set-TextValue $ws.Range("E29") "  +20.56%  "
# Row 31
Set-TextValue $ws.Range("E31") "  +8.90%  "
# Row 32
Set-TextValue $ws.Range("D32") "1.361"
Set-TextValue $ws.Range("E32") "  +2.97%  "
# Row 33
Set-TextValue $ws.Range("D33") "4.366"
Set-TextValue $ws.Range("E33") "  +5.77%  "
# Row 34
Set-TextValue $ws.Range("D34") "0.05224"
Set-TextValue $ws.Range("E34") "  +7.86%  "
# Row 35
Set-TextValue $ws.Range("D35") "1.220"
Set-TextValue $ws.Range("E35") "  +8.61%  "
# Row 36
Set-TextValue $ws.Range("D36") "0.7644"
Set-TextValue $ws.Range("E36") "  +10.40%  "
# Row 37
Set-TextValue $ws.Range("D37") "2.752"
Set-TextValue $ws.Range("E37") "  +1.34%  "
# Row 38
Set-TextValue $ws.Range("D38") "0.02014"
Set-TextValue $ws.Range("E38") "  +6.14%  "
# Row 39
Set-TextValue $ws.Range("D39") "2.921"
Set-TextValue $ws.Range("E39") "  +3.64%  "
# Row 40
Set-TextValue $ws.Range("D40") "6.651"
Set-TextValue $ws.Range("E40") "  +6.19%  "
# Row 41
Set-TextValue $ws.Range("D41") "0.4764"
Set-TextValue $ws.Range("E41") "  +12.41%  "
# Row 42
Set-TextValue $ws.Range("D42") "78.74"
Set-TextValue $ws.Range("E42") "  +3.12%  "
# Row 43
Set-TextValue $ws.Range("D43") "2.136"
Set-TextValue $ws.Range("E43") "  +10.08%  "
# Row 44
Set-TextValue $ws.Range("D44") "0.8622"
Set-TextValue $ws.Range("E44") "  +3.97%  "
# Row 45
Set-TextValue $ws.Range("D45") "104.68"
Set-TextValue $ws.Range("E45") "  +3.63%  "
# Row 46
Set-TextValue $ws.Range("D46") "0.9977"
Set-TextValue $ws.Range("E46") "  -0.12%  "
# Row 47
Set-TextValue $ws.Range("D47") "10.05"
Set-TextValue $ws.Range("E47") "  +3.00%  "
# Row 48
Set-TextValue $ws.Range("D48") "7.544"
Set-TextValue $ws.Range("E48") "  +7.99%  "
# Row 49
Set-TextValue $ws.Range("D49") "0.4402"
Set-TextValue $ws.Range("E49") "  +11.65%  "
# Row 50
Set-TextValue $ws.Range("D50") "36.90"
Set-TextValue $ws.Range("E50") "  +4.53%  "
# Row 51
Set-TextValue $ws.Range("D51") "0.1209"
Set-TextValue $ws.Range("E51") "  +14.17%  "
